$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.7466303911698731
$ws.Range("C2").Value = 0.09815769476978176
$ws.Range("D2").Value = 0.0599272770681285
$ws.Range("E2").Value = 0.4056186826292389
$ws.Range("F2").Value = 1.661456248505488
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("K2").Value = 0.767741738778966
$ws.Range("N2").Value = 2.233450337813224

$ws.Range("B3").Value = 0.6889620168639681
$ws.Range("C3").Value = 0.0865679908451682
$ws.Range("D3").Value = 0.06025035040060445
$ws.Range("E3").Value = 0.3537967305260565
$ws.Range("F3").Value = 1.624321692453606
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("K3").Value = 0.7010882620877226
$ws.Range("N3").Value = 2.231680028311615

$ws.Range("B4").Value = 0.6539997000339497
$ws.Range("C4").Value = 0.07948822867972183
$ws.Range("D4").Value = 0.06045457121038389
$ws.Range("E4").Value = 0.322097290657652
$ws.Range("F4").Value = 1.602408460569649
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("K4").Value = 0.6606010198522654
$ws.Range("N4").Value = 2.231203198771254

$ws.Range("B5").Value = 0.6398634919544008
$ws.Range("C5").Value = 0.07661200573544136
$ws.Range("D5").Value = 0.06053925492397561
$ws.Range("E5").Value = 0.3092071012994069
$ws.Range("F5").Value = 1.593700552284275
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("K5").Value = 0.644211038315575
$ws.Range("N5").Value = 2.231161623374504

$ws.Range("B6").Value = 0.6375228737719283
$ws.Range("C6").Value = 0.07613493600779009
$ws.Range("D6").Value = 0.06055340462221892
$ws.Range("E6").Value = 0.3070682905105286
$ws.Range("F6").Value = 1.592267976964834
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("K6").Value = 0.6414960338223352
$ws.Range("N6").Value = 2.23116392465478

$ws.Range("B7").Value = 0.6538086051822631
$ws.Range("C7").Value = 0.07944940354798291
$ws.Range("D7").Value = 0.06045570737814465
$ws.Range("E7").Value = 0.3219233405378503
$ws.Range("F7").Value = 1.602290125848327
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("K7").Value = 0.660379539616855
$ws.Range("N7").Value = 2.231202020588469

$ws.Range("B8").Value = 0.7266530778292122
$ws.Range("C8").Value = 0.09415382904560943
$ws.Range("D8").Value = 0.0600374522541518
$ws.Range("E8").Value = 0.3877241883698161
$ws.Range("F8").Value = 1.648467264238207
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("K8").Value = 0.7446678458734937
$ws.Range("N8").Value = 2.232712878586483

$ws.Range("B9").Value = 0.8730959334046133
$ws.Range("C9").Value = 0.1232930281187521
$ws.Range("D9").Value = 0.05926411483215155
$ws.Range("E9").Value = 0.5178241206303795
$ws.Range("F9").Value = 1.746126653108647
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("K9").Value = 0.9135028826585199
$ws.Range("N9").Value = 2.240549689117429

$ws.Range("B10").Value = 0.9829661082179086
$ws.Range("C10").Value = 0.1449103022625309
$ws.Range("D10").Value = 0.05872508360705453
$ws.Range("E10").Value = 0.6142289686525828
$ws.Range("F10").Value = 1.822307572166096
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("K10").Value = 1.039816409272447
$ws.Range("N10").Value = 2.249326500605306

$ws.Range("B11").Value = 1.033463459535255
$ws.Range("C11").Value = 0.1547951352413577
$ws.Range("D11").Value = 0.05848632285574773
$ws.Range("E11").Value = 0.6583021958963542
$ws.Range("F11").Value = 1.85794842826283
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("K11").Value = 1.097797327343415
$ws.Range("N11").Value = 2.253985464333553

$ws.Range("B12").Value = 1.052661201842398
$ws.Range("C12").Value = 0.1585459907333586
$ws.Range("D12").Value = 0.05839685001308581
$ws.Range("E12").Value = 0.675025915821152
$ws.Range("F12").Value = 1.87158808389924
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("K12").Value = 1.119829758211495
$ws.Range("N12").Value = 2.255846352783593

$ws.Range("B13").Value = 1.048523248194897
$ws.Range("C13").Value = 0.1577378295793039
$ws.Range("D13").Value = 0.05841607753247402
$ws.Range("E13").Value = 0.6714226002893469
$ws.Range("F13").Value = 1.868644145615434
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("K13").Value = 1.115081261622038
$ws.Range("N13").Value = 2.255441265881956

$ws.Range("B14").Value = 1.035041349615483
$ws.Range("C14").Value = 0.1551035643823582
$ws.Range("D14").Value = 0.05847894293854239
$ws.Range("E14").Value = 0.6596773660403557
$ws.Range("F14").Value = 1.859067691673545
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("K14").Value = 1.099608410502867
$ws.Range("N14").Value = 2.254136619201191

$ws.Range("B15").Value = 1.026793163721038
$ws.Range("C15").Value = 0.1534910119150652
$ws.Range("D15").Value = 0.05851757273013902
$ws.Range("E15").Value = 0.6524876006740357
$ws.Range("F15").Value = 1.853220535252348
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("K15").Value = 1.090140826052732
$ws.Range("N15").Value = 2.253350095344913

$ws.Range("B16").Value = 0.9796764873269694
$ws.Range("C16").Value = 0.1442653614284097
$ws.Range("D16").Value = 0.05874081806671327
$ws.Range("E16").Value = 0.6113533058149585
$ws.Range("F16").Value = 1.819998320978925
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("K16").Value = 1.036037816693721
$ws.Range("N16").Value = 2.249035504825969

$ws.Range("B17").Value = 0.9509051354525013
$ws.Range("C17").Value = 0.1386190373113152
$ws.Range("D17").Value = 0.05887943262303708
$ws.Range("E17").Value = 0.5861765870147622
$ws.Range("F17").Value = 1.799871097025104
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("K17").Value = 1.002981684275483
$ws.Range("N17").Value = 2.246559890815433

$ws.Range("B18").Value = 0.9344051926356087
$ws.Range("C18").Value = 0.135376202824574
$ws.Range("D18").Value = 0.05895976605164677
$ws.Range("E18").Value = 0.5717159718342941
$ws.Range("F18").Value = 1.788387145680247
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("K18").Value = 0.9840175490617753
$ws.Range("N18").Value = 2.245198653261212

$ws.Range("B19").Value = 0.9288269036002248
$ws.Range("C19").Value = 0.1342790456127148
$ws.Range("D19").Value = 0.05898706924660502
$ws.Range("E19").Value = 0.5668232704412475
$ws.Range("F19").Value = 1.784514752793427
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("K19").Value = 0.9776049682323276
$ws.Range("N19").Value = 2.244748500457916

$ws.Range("B20").Value = 0.9539628589747053
$ws.Range("C20").Value = 0.1392196007602138
$ws.Range("D20").Value = 0.05886461405390975
$ws.Range("E20").Value = 0.5888545637521077
$ws.Range("F20").Value = 1.802004070107586
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("K20").Value = 1.006495496531556
$ws.Range("N20").Value = 2.246816932184302

$ws.Range("B21").Value = 1.038999251794678
$ws.Range("C21").Value = 0.1558771006026518
$ws.Range("D21").Value = 0.05846045220395091
$ws.Range("E21").Value = 0.6631262770471551
$ws.Range("F21").Value = 1.861876628974812
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("K21").Value = 1.104151078642758
$ws.Range("N21").Value = 2.254517196181041

$ws.Range("B22").Value = 1.095016030675538
$ws.Range("C22").Value = 0.1668086914583569
$ws.Range("D22").Value = 0.05820179872321063
$ws.Range("E22").Value = 0.7118674783994976
$ws.Range("F22").Value = 1.901842408072923
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("K22").Value = 1.168420190436166
$ws.Range("N22").Value = 2.260113444487786

$ws.Range("B23").Value = 1.065078116239022
$ws.Range("C23").Value = 0.1609700666752758
$ws.Range("D23").Value = 0.05833934002416008
$ws.Range("E23").Value = 0.6858341230853
$ws.Range("F23").Value = 1.880434965373809
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("K23").Value = 1.13407728458003
$ws.Range("N23").Value = 2.257074770686316

$ws.Range("B24").Value = 0.9525803342702375
$ws.Range("C24").Value = 0.1389480757363799
$ws.Range("D24").Value = 0.05887131152912239
$ws.Range("E24").Value = 0.5876438076756614
$ws.Range("F24").Value = 1.801039480657337
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("K24").Value = 1.004906776677615
$ws.Range("N24").Value = 2.246700530616991

$ws.Range("B25").Value = 0.8330854769858433
$ws.Range("C25").Value = 0.1153751661605327
$ws.Range("D25").Value = 0.05946824915581317
$ws.Range("E25").Value = 0.4824970907171462
$ws.Range("F25").Value = 1.718936341768142
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("K25").Value = 0.8674377740411785
$ws.Range("N25").Value = 2.23790292703822
